$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FILESTREAM")

$ws.Range("A9").Value = "NHANES-2009-2010-UCPREG_F"
$ws.Range("B9").Value = "NHANES-LAB-RESULTS"
$ws.Range("C9").Value = "nhanes-kb:DPL-BECKMAN-COULTER-ICON-25-HCG-URINE-SERUM-TEST-KIT"
$ws.Range("E9").Value = "example@example.com"
$ws.Range("F9").Value = "Public"

$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:example@example.com")

$ws.Range("C3").Select()
Write-Output "done"
